$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Cenário Normal" table content (rows 7-8) ---
# Row 7: Actor input now holds the text that used to live on row 8 ("2. Categoriza o stock"
# is being replaced); System response is cleared (text moves to row 8 instead).
# D8 must be filled BEFORE C7 so that the shared-string table records
# "2. Processa o novo stock" ahead of "1.Regista o stock que chegou ao sistema",
# matching the authoring order captured in the workbook.
$ws.Range("D8").Value = "2. Processa o novo stock"
$ws.Range("C7").Value = "1.Regista o stock que chegou ao sistema"
$ws.Range("D7").Value = ""
$ws.Range("C8").Value = ""

# --- Remove the old "Cenário alternativo 1" block (rows 17-19) ---
$ws.Rows("17:19").Delete()

# --- Widen column C to fit the new wording ---
$ws.Columns("C").ColumnWidth = 40.666666666666664

# --- Zoom the sheet to 210% as recorded in the saved view state ---
$excel.ActiveWindow.Zoom = 210
